# Weekly fruit/vegetable price update: insert a new weekly record for
# "Terminal Hortofrutícola Agro Chillán - Cilantro" as row 59, pushing the
# existing rows 59-126 down to 60-127 (dimension grows from A1:R126 to A1:R127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 59 (shifts rows 59..126 down to 60..127,
# carrying the existing row styling with it).
$ws.Rows.Item(59).Insert()

# Populate the new row 59 with this week's record.
$ws.Cells.Item(59, 1).Value  = 7
$ws.Cells.Item(59, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(59, 3).Value  = "Ñuble"
$ws.Cells.Item(59, 4).Value  = 44897
$ws.Cells.Item(59, 5).Value  = 16
$ws.Cells.Item(59, 6).Value  = 100112040
$ws.Cells.Item(59, 7).Value  = "Cilantro"
$ws.Cells.Item(59, 8).Value  = "Sin especificar"
$ws.Cells.Item(59, 9).Value  = "Primera"
$ws.Cells.Item(59, 10).Value = 400
$ws.Cells.Item(59, 11).Value = 600
$ws.Cells.Item(59, 12).Value = 700
$ws.Cells.Item(59, 13).Value = 650
$ws.Cells.Item(59, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(59, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(59, 16).Value = 650
$ws.Cells.Item(59, 17).Value = 1
$ws.Cells.Item(59, 18).Value = "Hortaliza"
